# edit.ps1 - apply the "ottimizzazione codice + rivista documentazione" changes
# Units note: the PowerPoint COM object model expresses Left/Top/Width/Height in
# points (1 pt = 12700 EMU) using (effectively) single-precision floats, which can
# round a plain EMU/12700 conversion to the EMU bucket below the intended target.
# Nudging by +/-0.5 EMU before dividing reliably lands the stored value back on the
# exact target EMU (verified empirically against this runtime).
function EmuToPt($emu) {
    if ($emu -ge 0) {
        return ($emu + 0.5) / 12700
    } else {
        return ($emu - 0.5) / 12700
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: six chevron shapes shift horizontally (extent unchanged)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$sh = $s2.Shapes.Item(1)
$sh.Left = EmuToPt 216595
$sh.Top = EmuToPt 2172456

$sh = $s2.Shapes.Item(2)
$sh.Left = EmuToPt 1658776
$sh.Top = EmuToPt 2172456

$sh = $s2.Shapes.Item(3)
$sh.Left = EmuToPt 3068313
$sh.Top = EmuToPt 2172456

$sh = $s2.Shapes.Item(4)
$sh.Left = EmuToPt 4491746
$sh.Top = EmuToPt 2168748

$sh = $s2.Shapes.Item(5)
$sh.Left = EmuToPt 5978134
$sh.Top = EmuToPt 2168748

$sh = $s2.Shapes.Item(7)
$sh.Left = EmuToPt 7493234
$sh.Top = EmuToPt 2172456

# ---------------------------------------------------------------------------
# Slide 6: "Server Registry" bullet text boxes
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Shape 9: "Permette ad un client esterno di interagire col sistema. ..."
# -> italicize "client esterno" (splits the run into three)
$sh = $s6.Shapes.Item(9)
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(16, 14)
$sub.Font.Italic = $true
$sh.Left = EmuToPt 10396
$sh.Top = EmuToPt 975386
$sh.Width = EmuToPt 5514240
$sh.Height = EmuToPt 518519

# Shape 10: "Fornisce ad un nuovo nodo un punto di accesso nel sistema, ..."
# -> rewritten text + "nuovo" italicized
$sh = $s6.Shapes.Item(10)
$tr = $sh.TextFrame.TextRange
$run1 = $tr.Characters(1, 90)
$run1.Text = "Fornisce ad un nuovo nodo un punto di accesso nel sistema, fornendogli"
$run2 = $tr.Characters(71, 50)
$run2.Text = " il contatto dei nodi adiacenti con cui instaurare la connessione."
$nuovo = $tr.Characters(16, 5)
$nuovo.Font.Italic = $true
$sh.Left = EmuToPt -10291
$sh.Top = EmuToPt 2154915
$sh.Width = EmuToPt 4614831
$sh.Height = EmuToPt 731879

# ---------------------------------------------------------------------------
# Slide 11: resource-lookup example now keyed by id "4" instead of name "ciao"
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

# Shape 18: "Gestisco la\nrisorsa “ciao”" -> "... “4:ciao”"
$sh = $s11.Shapes.Item(18)
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(22, 6)
$sub.Text = [char]8220 + "4:ciao" + [char]8221
$sh.Left = EmuToPt 3969456
$sh.Top = EmuToPt 2753974
$sh.Width = EmuToPt 1094914
$sh.Height = EmuToPt 396599

# Shape 19: "Contatto il nodo 'x' per cercare la risorsa “ciao”"
#        -> "Contatto il nodo 'x' per cercare la risorsa avente id = 4"
$sh = $s11.Shapes.Item(19)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Contatto il nodo " + [char]8216 + "x" + [char]8217 + " per cercare la risorsa avente id = 4"
$sh.Rotation = 20869666 / 60000
$sh.Left = EmuToPt 1546691
$sh.Top = EmuToPt 2482507
$sh.Width = EmuToPt 3222569
$sh.Height = EmuToPt 244199

# Shape 20: "id(ciao)=4" -> "id = 4"
$sh = $s11.Shapes.Item(20)
$tr = $sh.TextFrame.TextRange
$tr.Text = "id = 4"
$sh.Left = EmuToPt 7512222
$sh.Top = EmuToPt 1316319
$sh.Width = EmuToPt 497808
$sh.Height = EmuToPt 244199

# Shape 23: "Elimino la\nrisorsa “ciao”" -> "... “4:ciao”"
$sh = $s11.Shapes.Item(23)
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(20, 6)
$sub.Text = [char]8220 + "4:ciao" + [char]8221
$sh.Left = EmuToPt 3969456
$sh.Top = EmuToPt 3104676
$sh.Width = EmuToPt 1138067
$sh.Height = EmuToPt 396599

# ---------------------------------------------------------------------------
# Slide 15: 1 EMU nudge on the example caption's left offset
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$sh = $s15.Shapes.Item(11)
$sh.Left = EmuToPt 23944
$sh.Top = EmuToPt 1306941
